$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell C1 ("note"), bold style like the other headers
$ws.Range("C1").Value = "note"
$ws.Range("C1").Font.Bold = $true

# Row 4: Erase Fail
$ws.Range("C4").Value = "Erase Fail"
$ws.Range("B4").Value = 1
$ws.Range("A4").Value = "PS: EH  EF  Recovery(reqId|VBA|msgType|cmd70)"

# Row 5: Prog Fail
$ws.Range("C5").Value = "Prog Fail"
$ws.Range("B5").Value = 1
$ws.Range("A5").Value = "PS: EH  PF Recovery(reqID|VBA|msgType Cmd70 blkType|PEC)"

# Row 6: UECC
$ws.Range("C6").Value = "UECC"
$ws.Range("B6").Value = 1
$ws.Range("A6").Value = "PS: EH  REH  UNRECOVERABLE(dgId|reqIdx|msgType|EHmap|VBA|savedMap|moreInfo)"

# Row 7: EPWR Failure
$ws.Range("A7").Value = "PS: BBM  EH  EPWR block bad(deVBA)"
$ws.Range("C7").Value = "EPWR Failure after recovery  fail (GBB)"
$ws.Range("B7").Value = 1

# New-row formatting: apply the Segoe UI / dark-navy font used for the new note column entries
$noteFont = $ws.Range("A4:A7")
$noteFont.Font.Name = "Segoe UI"
$noteFont.Font.Color = 5057303
$ws.Range("C5").Font.Name = "Segoe UI"
$ws.Range("C5").Font.Color = 5057303
$ws.Range("C7").Font.Name = "Segoe UI"
$ws.Range("C7").Font.Color = 5057303

# Row heights for the new rows
$ws.Rows.Item(4).RowHeight = 16.5
$ws.Rows.Item(5).RowHeight = 16.5
$ws.Rows.Item(6).RowHeight = 16.5
$ws.Rows.Item(7).RowHeight = 16.5

# Column A widened to fit the longer note text
$ws.Columns.Item(1).ColumnWidth = 59.6

# Selection / active cell as left by the author
$ws.Range("A6").Select()
